# "new moose transect transform"
#
# 1) event sheet: drop the eventRemarks column (C) - that remark moves
#    to the occurrence sheet's occurrenceRemarks instead (out of scope here).
# 2) measurementOrFact sheet: drop the last "observation-comments" measurement
#    row entirely, and re-order the per-occurrence measurement rows so the
#    "antler-configuration" row (previously first) becomes the last of that
#    block instead (with its occurrenceID column), while the rest of the
#    per-occurrence rows lose the occurrenceID column value.

$wb = $excel.ActiveWorkbook

# --- Sheet "event": remove column C (eventRemarks) ---------------------
$wsEvent = $wb.Worksheets.Item("event")
$wsEvent.Columns.Item(3).Delete()

# --- Sheet "measurementOrFact": reorder rows 5-15 -----------------------
$wsMof = $wb.Worksheets.Item("measurementOrFact")

# Row 5 ("antler-configuration") is moving down to become the new last
# data row (14), after row 15 ("observation-comments") is dropped. Grab
# its values before the shift destroys them.
$bVal = $wsMof.Cells.Item(5, 2).Text
$cVal = $wsMof.Cells.Item(5, 3).Text
$dVal = $wsMof.Cells.Item(5, 4).Text
$eVal = $wsMof.Cells.Item(5, 5).Text
$fVal = $wsMof.Cells.Item(5, 6).Text

# Deleting row 5 shifts old rows 6..15 up to 5..14, which both re-orders
# the remaining measurements and drops the old "observation-comments" row
# (15) for free, since it no longer exists once everything shifts up.
$wsMof.Rows.Item(5).Delete()

# Only the final measurement row of each occurrence block carries the
# occurrenceID (column F); clear it from the rows that shifted into
# positions 5-13 (it was only valid there because of the old row 5 row).
$wsMof.Range("F5:F13").ClearContents()

# Restore the antler-configuration measurement as the new final row (14),
# including its occurrenceID.
$wsMof.Cells.Item(14, 1).Value = "MU 749:Lowland:1:0"
$wsMof.Cells.Item(14, 2).Value = $bVal
$wsMof.Cells.Item(14, 3).Value = $cVal
if ($dVal -ne "") {
    $wsMof.Cells.Item(14, 4).Value = $dVal
}
$wsMof.Cells.Item(14, 5).Value = $eVal
$wsMof.Cells.Item(14, 6).Value = $fVal
